# Change demand growth and availability factor for LDVs
#
# - UCT1 sheet view: was the active/selected tab, scrolled to A4 -> no longer
#   the selected tab, scrolled further down (A25). Selected cell stays F76.
# - Sheet9 sheet view: becomes the active/selected tab; selection moves from
#   E8 to D15:D17 (the three growth-factor / availability-factor cells that
#   are being edited).
# - Sheet9 D15:D17 (TRAPS/TRAPM/TRAPL demand growth & availability factors
#   for LDVs): 1.1 -> 1.5. Downstream formulas in rows 5, 7 and 9
#   (G:AG) recompute automatically from this input change.

$wb = $excel.ActiveWorkbook

# --- UCT1: leave the previously-active tab, scroll down, keep selection ---
$uct1 = $wb.Worksheets.Item("UCT1")
[void]$uct1.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
[void]$uct1.Range("F76").Select()

# --- Sheet9: the real edit - bump the LDV demand growth / availability ---
$sheet9 = $wb.Worksheets.Item("Sheet9")
[void]$sheet9.Activate()

$sheet9.Range("D15").Value = 1.5
$sheet9.Range("D16").Value = 1.5
$sheet9.Range("D17").Value = 1.5

# Sheet9 becomes the active tab, with D15:D17 selected (the edited cells)
[void]$sheet9.Range("D15:D17").Select()
